{"js": "// This script replaces the date line and each \"a\\u00f7b=\" expression in the\n// table with the new values specified by the commit's diff. Every source\n// string in this document is unique, so a body-wide search-and-replace for\n// each exact string is sufficient and safe (no accidental double edits).\nconst replacements = [\n  [\"2025-07-29 Tuesday\", \"2025-07-30 Wednesday\"],\n  [\"240\u00f75=\", \"354\u00f77=\"],\n  [\"362\u00f73=\", \"416\u00f76=\"],\n  [\"866\u00f74=\", \"887\u00f72=\"],\n  [\"794\u00f74=\", \"444\u00f73=\"],\n  [\"325\u00f75=\", \"898\u00f72=\"],\n  [\"671\u00f75=\", \"861\u00f74=\"],\n  [\"664\u00f77=\", \"539\u00f77=\"],\n  [\"848\u00f73=\", \"994\u00f75=\"],\n  [\"754\u00f73=\", \"901\u00f78=\"],\n  [\"909\u00f75=\", \"760\u00f79=\"],\n  [\"110\u00f72=\", \"359\u00f78=\"],\n  [\"926\u00f77=\", \"461\u00f77=\"],\n  [\"603\u00f77=\", \"192\u00f78=\"],\n  [\"622\u00f78=\", \"541\u00f73=\"],\n  [\"902\u00f76=\", \"966\u00f75=\"],\n  [\"128\u00f76=\", \"620\u00f76=\"],\n  [\"790\u00f75=\", \"330\u00f77=\"],\n  [\"610\u00f77=\", \"513\u00f77=\"],\n  [\"589\u00f79=\", \"237\u00f72=\"],\n  [\"829\u00f77=\", \"223\u00f75=\"],\n  [\"666\u00f76=\", \"740\u00f73=\"],\n  [\"296\u00f74=\", \"401\u00f73=\"],\n  [\"756\u00f77=\", \"118\u00f74=\"],\n  [\"921\u00f78=\", \"943\u00f72=\"],\n  [\"379\u00f76=\", \"696\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Updates the worksheet date header and every \"a(div)b=\" expression inside\n# the practice-problem table to the new values from the commit's diff.\n# Each original string occurs exactly once in the document, so a simple\n# Find/Replace (scoped to the whole story with wdReplaceAll) for every\n# old->new pair reproduces the edit deterministically.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '2025-07-29 Tuesday'; New = '2025-07-30 Wednesday' },\n    @{ Old = '240\u00f75='; New = '354\u00f77=' },\n    @{ Old = '362\u00f73='; New = '416\u00f76=' },\n    @{ Old = '866\u00f74='; New = '887\u00f72=' },\n    @{ Old = '794\u00f74='; New = '444\u00f73=' },\n    @{ Old = '325\u00f75='; New = '898\u00f72=' },\n    @{ Old = '671\u00f75='; New = '861\u00f74=' },\n    @{ Old = '664\u00f77='; New = '539\u00f77=' },\n    @{ Old = '848\u00f73='; New = '994\u00f75=' },\n    @{ Old = '754\u00f73='; New = '901\u00f78=' },\n    @{ Old = '909\u00f75='; New = '760\u00f79=' },\n    @{ Old = '110\u00f72='; New = '359\u00f78=' },\n    @{ Old = '926\u00f77='; New = '461\u00f77=' },\n    @{ Old = '603\u00f77='; New = '192\u00f78=' },\n    @{ Old = '622\u00f78='; New = '541\u00f73=' },\n    @{ Old = '902\u00f76='; New = '966\u00f75=' },\n    @{ Old = '128\u00f76='; New = '620\u00f76=' },\n    @{ Old = '790\u00f75='; New = '330\u00f77=' },\n    @{ Old = '610\u00f77='; New = '513\u00f77=' },\n    @{ Old = '589\u00f79='; New = '237\u00f72=' },\n    @{ Old = '829\u00f77='; New = '223\u00f75=' },\n    @{ Old = '666\u00f76='; New = '740\u00f73=' },\n    @{ Old = '296\u00f74='; New = '401\u00f73=' },\n    @{ Old = '756\u00f77='; New = '118\u00f74=' },\n    @{ Old = '921\u00f78='; New = '943\u00f72=' },\n    @{ Old = '379\u00f76='; New = '696\u00f79=' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
